$d = $word.ActiveDocument

# Update the date/weekday heading
$d.Content.Find.Execute("2025-09-17 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-09-18 Thursday", 2)

# Update the 25 division problems laid out in the 5x5 table (rows 1,5,9,13,17
# of the underlying w:tbl correspond to the data rows; other rows are blank
# spacer rows). Addressing by Table.Cell(row, col) avoids ambiguity from the
# duplicated "65÷6=" values that appear twice in the source document.
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text  = "78÷8="
$t.Cell(1,2).Range.Text  = "75÷9="
$t.Cell(1,3).Range.Text  = "43÷2="
$t.Cell(1,4).Range.Text  = "22÷5="
$t.Cell(1,5).Range.Text  = "45÷5="

$t.Cell(5,1).Range.Text  = "33÷5="
$t.Cell(5,2).Range.Text  = "86÷3="
$t.Cell(5,3).Range.Text  = "68÷2="
$t.Cell(5,4).Range.Text  = "98÷2="
$t.Cell(5,5).Range.Text  = "57÷7="

$t.Cell(9,1).Range.Text  = "35÷8="
$t.Cell(9,2).Range.Text  = "94÷2="
$t.Cell(9,3).Range.Text  = "98÷7="
$t.Cell(9,4).Range.Text  = "89÷7="
$t.Cell(9,5).Range.Text  = "24÷5="

$t.Cell(13,1).Range.Text = "52÷5="
$t.Cell(13,2).Range.Text = "34÷5="
$t.Cell(13,3).Range.Text = "48÷8="
$t.Cell(13,4).Range.Text = "44÷5="
$t.Cell(13,5).Range.Text = "56÷4="

$t.Cell(17,1).Range.Text = "30÷4="
$t.Cell(17,2).Range.Text = "74÷5="
$t.Cell(17,3).Range.Text = "70÷6="
$t.Cell(17,4).Range.Text = "60÷5="
$t.Cell(17,5).Range.Text = "20÷2="
